$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "78.852.33"
$ws.Range("E2").Value = "  +3.04%  "
$ws.Range("D3").Value = "3.175.68"
$ws.Range("E3").Value = "  +4.55%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'205.61"
$ws.Range("E5").Value = "  +2.28%  "
$ws.Range("D6").Value = "'629.90"
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'0.226"
$ws.Range("E8").Value = "  +10.62%  "
$ws.Range("D9").Value = "'0.580"
$ws.Range("E9").Value = "  +4.98%  "
$ws.Range("D10").Value = "3.172.72"
$ws.Range("E10").Value = "  +4.56%  "
$ws.Range("D11").Value = "'0.574"
$ws.Range("E11").Value = "  +31.35%  "
$ws.Range("E12").Value = "  +2.56%  "
$ws.Range("D13").Value = "'5.41"
$ws.Range("E13").Value = "  +6.66%  "
$ws.Range("D14").Value = "3.755.59"
$ws.Range("E14").Value = "  +4.41%  "
$ws.Range("E15").Value = "  +16.64%  "
$ws.Range("D16").Value = "'31.42"
$ws.Range("E16").Value = "  +6.82%  "
$ws.Range("D17").Value = "78.782.30"
$ws.Range("E17").Value = "  +3.09%  "
$ws.Range("D18").Value = "3.153.20"
$ws.Range("E18").Value = "  +4.43%  "
$ws.Range("E19").Value = "  +6.97%  "
$ws.Range("D20").Value = "'9.30"
$ws.Range("E20").Value = "  +2.33%  "
$ws.Range("D21").Value = "'426.86"
$ws.Range("E21").Value = "  +13.47%  "
$ws.Range("D22").Value = "'2.82"
$ws.Range("E22").Value = "  +23.44%  "
$ws.Range("D23").Value = "'4.92"
$ws.Range("E23").Value = "  +12.94%  "
$ws.Range("D24").Value = "'6.81"
$ws.Range("E24").Value = "  +5.32%  "
$ws.Range("D25").Value = "3.332.15"
$ws.Range("E25").Value = "  +3.92%  "
$ws.Range("D26").Value = "'4.72"
$ws.Range("E26").Value = "  +7.09%  "
$ws.Range("D27").Value = "'75.63"
$ws.Range("E27").Value = "  +3.50%  "
$ws.Range("D28").Value = "'10.88"
$ws.Range("E28").Value = "  +10.09%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.30%  "
$ws.Range("E30").Value = "  +3.55%  "
$ws.Range("D31").Value = "'0.999"
$ws.Range("E31").Value = "  +0.19%  "
$ws.Range("E32").Value = "  +5.75%  "
$ws.Range("E33").Value = "  +4.56%  "
$ws.Range("D34").Value = "'510.26"
$ws.Range("E34").Value = "  -0.86%  "
$ws.Range("E35").Value = "  +0.51%  "
$ws.Range("B36").Value = "EthereumClassic"
$ws.Range("C36").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D36").Value = "'22.79"
$ws.Range("E36").Value = "  +8.80%  "
$ws.Range("B37").Value = "Cronos"
$ws.Range("C37").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D37").Value = "'0.124"
$ws.Range("E37").Value = "  +18.58%  "
$ws.Range("E38").Value = "  +19.12%  "
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("B40").Value = "PolygonEcosystemToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D40").Value = "'0.396"
$ws.Range("E40").Value = "  +3.66%  "
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").Value = "'163.98"
$ws.Range("E41").Value = "  -0.15%  "
$ws.Range("D42").Value = "'19.96"
$ws.Range("E42").Value = "  -0.23%  "
$ws.Range("E43").Value = "  -0.15%  "
$ws.Range("D44").Value = "'191.14"
$ws.Range("E44").Value = "  -1.04%  "
$ws.Range("D45").Value = "'5.38"
$ws.Range("E45").Value = "  +6.07%  "
$ws.Range("D46").Value = "'0.803"
$ws.Range("E46").Value = "  +14.52%  "
$ws.Range("D47").Value = "'1.78"
$ws.Range("E47").Value = "  +6.99%  "
$ws.Range("E48").Value = "  +3.26%  "
$ws.Range("D49").Value = "'42.53"
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("D50").Value = "'2.48"
$ws.Range("E50").Value = "  +4.86%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "'24.90"
$ws.Range("E51").Value = "  +10.63%  "
